$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2798.340480058071
